# Remove strikethrough formatting from task 2) text, and apply
# strikethrough formatting to the first sentence of task 3) text
# (text content itself is unchanged; only the <w:strike/> run
# formatting moves from task 2) to the start of task 3)).

$d = $word.ActiveDocument

function Unstrike-Text($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text to unstrike: $searchText"
    }
    $rng.Font.StrikeThrough = 0
}

function Strike-Text($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text to strike: $searchText"
    }
    $rng.Font.StrikeThrough = 1
}

# --- Task 2) : remove all strikethrough ---
Unstrike-Text("Redactor text-deschidera editarea afisarea fisier text")
Unstrike-Text("cautarea unui subsir in text(cautarea trebuie facuta de la pozitia curenta pana la pozitia dorita), inlocuirea unui subsir cu alt subsir(la solicitare sau toate),")
Unstrike-Text("permite concomitent lucru cu mai multe fisiere")
Unstrike-Text("modificarea fontului pe subsirurui")

# --- Task 3) : add strikethrough to the first sentence only ---
Strike-Text("Aplicatie cu baza de date-lucru cu baza de date cu minim trei tabele")

Write-Output "Done"
